# Actualización automática 2025-10-23 17:30:09
# Updates sales figures for GUERRERO FAREZ FABIAN MAURICIO across the three sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M4").Value = 5188.32
$wsVentasGrupo.Range("D7").Value = 285.12
$wsVentasGrupo.Range("K16").Value = 1045.44
$wsVentasGrupo.Range("D56").Value = "8 de 54"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F4").Value = 6808.31
$wsVentaMensual.Range("F7").Value = 285.12
$wsVentaMensual.Range("F16").Value = 4800.23
$wsVentaMensual.Range("F60").Value = 64473.37

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D3").Value = 5408.64
$wsCumplimiento.Range("E3").Value = 12260.5070988183
$wsCumplimiento.Range("F3").Value = 0.3061064560587493

$wsCumplimiento.Range("D10").Value = 8998.25
$wsCumplimiento.Range("E10").Value = -5117.17016465608
$wsCumplimiento.Range("F10").Value = 2.31849134306783

$wsCumplimiento.Range("D12").Value = 32605.91
$wsCumplimiento.Range("E12").Value = 20057.21
$wsCumplimiento.Range("F12").Value = 0.6191412510310821

$wsCumplimiento.Range("D14").Value = 62378.88
$wsCumplimiento.Range("E14").Value = 36637.62661190614
$wsCumplimiento.Range("F14").Value = 0.6299846574520467
